# Commit: "Move a slide to rec"
#
# The last two slides of the deck -- "Move to recitation" and
# "Testing data abstractions" -- are being pulled out of the lecture
# deck (they're moving into the recitation materials instead), so we
# just delete them from this presentation.

$p = $ppt.ActivePresentation

function Get-SlideTitle($slide) {
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            return $shp.TextFrame.TextRange.Text
        }
    }
    return ""
}

$titlesToRemove = @("Move to recitation", "Testing data abstractions")

# Walk backwards so deleting a slide doesn't shift the index of the
# slides we still need to examine/remove.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $title = Get-SlideTitle $slide
    if ($titlesToRemove -contains $title) {
        $slide.Delete()
    }
}
